$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value2 = "ECs"
$ws.Range("B2").Value2 = "Ntn1"
$ws.Range("C2").Value2 = "Neo1"
$ws.Range("D2").Value2 = "ECs"
$ws.Range("E2").Value2 = 2
$ws.Range("F2").Value2 = 0.6666666666666666
$ws.Range("G2").Value2 = 0.9305633333333333
$ws.Range("H2").Value2 = 2.79169
$ws.Range("I2").Value2 = 0.01768777137856805
$ws.Range("J2").Value2 = 0.01768777137856806
$ws.Range("K2").Value2 = 3
$ws.Range("L2").Value2 = 1
$ws.Range("M2").Value2 = 2.442429333333334
$ws.Range("N2").Value2 = 7.327288
$ws.Range("O2").Value2 = 0.08913295894744963
$ws.Range("P2").Value2 = 0.08913295894744963
$ws.Range("Q2").Value2 = 2.272835181857778
$ws.Range("R2").Value2 = 20.45551663672
$ws.Range("S2").Value2 = 0.001576563400157781
$ws.Range("T2").Value2 = 0.001576563400157781

# Row 3
$ws.Range("A3").Value2 = "ECs"
$ws.Range("B3").Value2 = "Ntn1"
$ws.Range("C3").Value2 = "Neo1"
$ws.Range("D3").Value2 = "FAPs"
$ws.Range("E3").Value2 = 2
$ws.Range("F3").Value2 = 0.6666666666666666
$ws.Range("G3").Value2 = 0.9305633333333333
$ws.Range("H3").Value2 = 2.79169
$ws.Range("I3").Value2 = 0.01768777137856805
$ws.Range("J3").Value2 = 0.01768777137856806
$ws.Range("K3").Value2 = 3
$ws.Range("L3").Value2 = 1
$ws.Range("M3").Value2 = 15.82990933333333
$ws.Range("N3").Value2 = 47.489728
$ws.Range("O3").Value2 = 0.5776898596383203
$ws.Range("P3").Value2 = 0.5776898596383203
$ws.Range("Q3").Value2 = 14.73073319559111
$ws.Range("R3").Value2 = 132.57659876032
$ws.Range("S3").Value2 = 0.01021804616499968
$ws.Range("T3").Value2 = 0.01021804616499968

# Row 4
$ws.Range("A4").Value2 = "ECs"
$ws.Range("B4").Value2 = "Ntn1"
$ws.Range("C4").Value2 = "Neo1"
$ws.Range("D4").Value2 = "sCs"
$ws.Range("E4").Value2 = 2
$ws.Range("F4").Value2 = 0.6666666666666666
$ws.Range("G4").Value2 = 0.9305633333333333
$ws.Range("H4").Value2 = 2.79169
$ws.Range("I4").Value2 = 0.01768777137856805
$ws.Range("J4").Value2 = 0.01768777137856806
$ws.Range("K4").Value2 = 3
$ws.Range("L4").Value2 = 1
$ws.Range("M4").Value2 = 9.129750999999999
$ws.Range("N4").Value2 = 27.389253
$ws.Range("O4").Value2 = 0.3331771814142301
$ws.Range("P4").Value2 = 0.3331771814142301
$ws.Range("Q4").Value2 = 8.495811523063331
$ws.Range("R4").Value2 = 76.46230370756999
$ws.Range("S4").Value2 = 0.005893161813410595
$ws.Range("T4").Value2 = 0.005893161813410596

# Row 5
$ws.Range("A5").Value2 = "FAPs"
$ws.Range("B5").Value2 = "Ntn1"
$ws.Range("C5").Value2 = "Neo1"
$ws.Range("D5").Value2 = "ECs"
$ws.Range("E5").Value2 = 3
$ws.Range("F5").Value2 = 1
$ws.Range("G5").Value2 = 44.154177
$ws.Range("H5").Value2 = 132.462531
$ws.Range("I5").Value2 = 0.8392647337471152
$ws.Range("J5").Value2 = 0.8392647337471153
$ws.Range("K5").Value2 = 3
$ws.Range("L5").Value2 = 1
$ws.Range("M5").Value2 = 2.442429333333334
$ws.Range("N5").Value2 = 7.327288
$ws.Range("O5").Value2 = 0.08913295894744963
$ws.Range("P5").Value2 = 0.08913295894744963
$ws.Range("Q5").Value2 = 107.843457093992
$ws.Range("R5").Value2 = 970.5911138459281
$ws.Range("S5").Value2 = 0.07480614905912386
$ws.Range("T5").Value2 = 0.07480614905912387

# Row 6
$ws.Range("A6").Value2 = "FAPs"
$ws.Range("B6").Value2 = "Ntn1"
$ws.Range("C6").Value2 = "Neo1"
$ws.Range("D6").Value2 = "FAPs"
$ws.Range("E6").Value2 = 3
$ws.Range("F6").Value2 = 1
$ws.Range("G6").Value2 = 44.154177
$ws.Range("H6").Value2 = 132.462531
$ws.Range("I6").Value2 = 0.8392647337471152
$ws.Range("J6").Value2 = 0.8392647337471153
$ws.Range("K6").Value2 = 3
$ws.Range("L6").Value2 = 1
$ws.Range("M6").Value2 = 15.82990933333333
$ws.Range("N6").Value2 = 47.489728
$ws.Range("O6").Value2 = 0.5776898596383203
$ws.Range("P6").Value2 = 0.5776898596383203
$ws.Range("Q6").Value2 = 698.9566185979521
$ws.Range("R6").Value2 = 6290.609567381569
$ws.Range("S6").Value2 = 0.4848347262377632
$ws.Range("T6").Value2 = 0.4848347262377633

# Row 7
$ws.Range("A7").Value2 = "FAPs"
$ws.Range("B7").Value2 = "Ntn1"
$ws.Range("C7").Value2 = "Neo1"
$ws.Range("D7").Value2 = "sCs"
$ws.Range("E7").Value2 = 3
$ws.Range("F7").Value2 = 1
$ws.Range("G7").Value2 = 44.154177
$ws.Range("H7").Value2 = 132.462531
$ws.Range("I7").Value2 = 0.8392647337471152
$ws.Range("J7").Value2 = 0.8392647337471153
$ws.Range("K7").Value2 = 3
$ws.Range("L7").Value2 = 1
$ws.Range("M7").Value2 = 9.129750999999999
$ws.Range("N7").Value2 = 27.389253
$ws.Range("O7").Value2 = 0.3331771814142301
$ws.Range("P7").Value2 = 0.3331771814142301
$ws.Range("Q7").Value2 = 403.116641619927
$ws.Range("R7").Value2 = 3628.049774579343
$ws.Range("S7").Value2 = 0.2796238584502281
$ws.Range("T7").Value2 = 0.2796238584502281

# Row 8
$ws.Range("A8").Value2 = "sCs"
$ws.Range("B8").Value2 = "Ntn1"
$ws.Range("C8").Value2 = "Neo1"
$ws.Range("D8").Value2 = "ECs"
$ws.Range("E8").Value2 = 3
$ws.Range("F8").Value2 = 1
$ws.Range("G8").Value2 = 7.525807
$ws.Range("H8").Value2 = 22.577421
$ws.Range("I8").Value2 = 0.1430474948743168
$ws.Range("J8").Value2 = 0.1430474948743168
$ws.Range("K8").Value2 = 3
$ws.Range("L8").Value2 = 1
$ws.Range("M8").Value2 = 2.442429333333334
$ws.Range("N8").Value2 = 7.327288
$ws.Range("O8").Value2 = 0.08913295894744963
$ws.Range("P8").Value2 = 0.08913295894744963
$ws.Range("Q8").Value2 = 18.38125177380534
$ws.Range("R8").Value2 = 165.431265964248
$ws.Range("S8").Value2 = 0.01275024648816799
$ws.Range("T8").Value2 = 0.01275024648816799

# Row 9
$ws.Range("A9").Value2 = "sCs"
$ws.Range("B9").Value2 = "Ntn1"
$ws.Range("C9").Value2 = "Neo1"
$ws.Range("D9").Value2 = "FAPs"
$ws.Range("E9").Value2 = 3
$ws.Range("F9").Value2 = 1
$ws.Range("G9").Value2 = 7.525807
$ws.Range("H9").Value2 = 22.577421
$ws.Range("I9").Value2 = 0.1430474948743168
$ws.Range("J9").Value2 = 0.1430474948743168
$ws.Range("K9").Value2 = 3
$ws.Range("L9").Value2 = 1
$ws.Range("M9").Value2 = 15.82990933333333
$ws.Range("N9").Value2 = 47.489728
$ws.Range("O9").Value2 = 0.5776898596383203
$ws.Range("P9").Value2 = 0.5776898596383203
$ws.Range("Q9").Value2 = 119.1328424701653
$ws.Range("R9").Value2 = 1072.195582231488
$ws.Range("S9").Value2 = 0.0826370872355574
$ws.Range("T9").Value2 = 0.0826370872355574

# Row 10
$ws.Range("A10").Value2 = "sCs"
$ws.Range("B10").Value2 = "Ntn1"
$ws.Range("C10").Value2 = "Neo1"
$ws.Range("D10").Value2 = "sCs"
$ws.Range("E10").Value2 = 3
$ws.Range("F10").Value2 = 1
$ws.Range("G10").Value2 = 7.525807
$ws.Range("H10").Value2 = 22.577421
$ws.Range("I10").Value2 = 0.1430474948743168
$ws.Range("J10").Value2 = 0.1430474948743168
$ws.Range("K10").Value2 = 3
$ws.Range("L10").Value2 = 1
$ws.Range("M10").Value2 = 9.129750999999999
$ws.Range("N10").Value2 = 27.389253
$ws.Range("O10").Value2 = 0.3331771814142301
$ws.Range("P10").Value2 = 0.3331771814142301
$ws.Range("Q10").Value2 = 68.708743984057
$ws.Range("R10").Value2 = 618.378695856513
$ws.Range("S10").Value2 = 0.04766016115059139
$ws.Range("T10").Value2 = 0.04766016115059139

